# Applies the Contra Costa County testing-data revision:
#  - corrects historical cumulative test counts (column E) and several
#    related derived/base cells (B, C, D, F, G) for existing rows
#  - appends one new day of data as row 272 (2020-12-27)
# Values below are taken verbatim from the target OOXML diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows ---
$ws.Range("E2").Value = 5779
$ws.Range("E3").Value = 6151
$ws.Range("E4").Value = 6587
$ws.Range("E5").Value = 6776
$ws.Range("E6").Value = 6951
$ws.Range("E7").Value = 7328
$ws.Range("E8").Value = 7717
$ws.Range("E9").Value = 8061
$ws.Range("E10").Value = 8385
$ws.Range("E11").Value = 8864
$ws.Range("E12").Value = 9047
$ws.Range("E13").Value = 9163
$ws.Range("E14").Value = 9522
$ws.Range("E15").Value = 9861
$ws.Range("E16").Value = 10214
$ws.Range("E17").Value = 10599
$ws.Range("E18").Value = 10949
$ws.Range("E19").Value = 11106
$ws.Range("E20").Value = 11252
$ws.Range("E21").Value = 11620
$ws.Range("E22").Value = 12029
$ws.Range("E23").Value = 12549
$ws.Range("E24").Value = 13108
$ws.Range("E25").Value = 13677
$ws.Range("E26").Value = 13909
$ws.Range("E27").Value = 14138
$ws.Range("E28").Value = 14748
$ws.Range("E29").Value = 15289
$ws.Range("E30").Value = 15982
$ws.Range("E31").Value = 16830
$ws.Range("E32").Value = 17559
$ws.Range("E33").Value = 17848
$ws.Range("E34").Value = 18110
$ws.Range("E35").Value = 18892
$ws.Range("E36").Value = 19570
$ws.Range("E37").Value = 20541
$ws.Range("E38").Value = 21438
$ws.Range("E39").Value = 22581
$ws.Range("F40").Value = 383
$ws.Range("G40").Value = 0.020328381548
$ws.Range("G41").Value = 0.020596727622
$ws.Range("G42").Value = 0.016678752719
$ws.Range("G43").Value = 0.014323345406
$ws.Range("G44").Value = 0.014532536735
$ws.Range("G45").Value = 0.015911673973
$ws.Range("G46").Value = 0.018534119629
$ws.Range("E55").Value = 36590
$ws.Range("F55").Value = 379
$ws.Range("G55").Value = 0.026626406807
$ws.Range("E56").Value = 37041
$ws.Range("G56").Value = 0.02573356347
$ws.Range("E57").Value = 38259
$ws.Range("G57").Value = 0.028123947457
$ws.Range("E58").Value = 39363
$ws.Range("G58").Value = 0.028324287652
$ws.Range("E59").Value = 40510
$ws.Range("G59").Value = 0.031276415891
$ws.Range("E60").Value = 41918
$ws.Range("G60").Value = 0.029860330711
$ws.Range("E61").Value = 42655
$ws.Range("G61").Value = 0.028553693358
$ws.Range("E62").Value = 43052
$ws.Range("E63").Value = 44424
$ws.Range("E64").Value = 45643
$ws.Range("E65").Value = 46814
$ws.Range("E66").Value = 48165
$ws.Range("E67").Value = 49800
$ws.Range("E68").Value = 50409
$ws.Range("E69").Value = 50897
$ws.Range("E70").Value = 52730
$ws.Range("E71").Value = 54578
$ws.Range("E72").Value = 56577
$ws.Range("E73").Value = 58736
$ws.Range("E74").Value = 61004
$ws.Range("E75").Value = 61849
$ws.Range("E76").Value = 62441
$ws.Range("E77").Value = 64668
$ws.Range("E78").Value = 67213
$ws.Range("E79").Value = 69796
$ws.Range("E80").Value = 72332
$ws.Range("E81").Value = 74783
$ws.Range("E82").Value = 75743
$ws.Range("E83").Value = 76445
$ws.Range("F83").Value = 702
$ws.Range("G83").Value = 0.036846615252
$ws.Range("E84").Value = 79067
$ws.Range("G84").Value = 0.039516633099
$ws.Range("E85").Value = 81862
$ws.Range("G85").Value = 0.042255444057
$ws.Range("E86").Value = 84906
$ws.Range("G86").Value = 0.044672402382
$ws.Range("E87").Value = 88291
$ws.Range("G87").Value = 0.046243498966
$ws.Range("E88").Value = 91281
$ws.Range("G88").Value = 0.048733179779
$ws.Range("E89").Value = 92364
$ws.Range("G89").Value = 0.050057156609
$ws.Range("E90").Value = 93164
$ws.Range("E91").Value = 96862
$ws.Range("E92").Value = 100358
$ws.Range("E93").Value = 103910
$ws.Range("E94").Value = 107852
$ws.Range("E95").Value = 110513
$ws.Range("E96").Value = 111641
$ws.Range("E97").Value = 112781
$ws.Range("E98").Value = 116877
$ws.Range("E99").Value = 120995
$ws.Range("E100").Value = 125271
$ws.Range("E101").Value = 129897
$ws.Range("E102").Value = 134378
$ws.Range("E103").Value = 136457
$ws.Range("E104").Value = 137794
$ws.Range("E105").Value = 142434
$ws.Range("E106").Value = 146224
$ws.Range("E107").Value = 150587
$ws.Range("E108").Value = 154799
$ws.Range("E109").Value = 159028
$ws.Range("E110").Value = 161270
$ws.Range("E111").Value = 162544
$ws.Range("E112").Value = 166820
$ws.Range("E113").Value = 170679
$ws.Range("E114").Value = 174373
$ws.Range("E115").Value = 177739
$ws.Range("E116").Value = 181441
$ws.Range("E117").Value = 183091
$ws.Range("F117").Value = 1650
$ws.Range("G117").Value = 0.06644975023999999
$ws.Range("E118").Value = 184355
$ws.Range("G118").Value = 0.06670945853
$ws.Range("E119").Value = 188561
$ws.Range("G119").Value = 0.06444045812
$ws.Range("E120").Value = 192518
$ws.Range("G120").Value = 0.065021292183
$ws.Range("E121").Value = 196408
$ws.Range("G121").Value = 0.06466984343
$ws.Range("E122").Value = 200116
$ws.Range("G122").Value = 0.063055816239
$ws.Range("E123").Value = 203685
$ws.Range("G123").Value = 0.064242042798
$ws.Range("E124").Value = 205365
$ws.Range("E125").Value = 206527
$ws.Range("E126").Value = 210420
$ws.Range("E127").Value = 213786
$ws.Range("E128").Value = 217099
$ws.Range("E129").Value = 220806
$ws.Range("E130").Value = 224170
$ws.Range("E131").Value = 225893
$ws.Range("E132").Value = 226970
$ws.Range("E133").Value = 230899
$ws.Range("E134").Value = 234338
$ws.Range("E135").Value = 238226
$ws.Range("E136").Value = 241526
$ws.Range("E137").Value = 244920
$ws.Range("E138").Value = 246494
$ws.Range("E139").Value = 247280
$ws.Range("E140").Value = 250958
$ws.Range("E141").Value = 254069
$ws.Range("E142").Value = 257162
$ws.Range("E143").Value = 260352
$ws.Range("E144").Value = 263408
$ws.Range("E145").Value = 264749
$ws.Range("F145").Value = 1341
$ws.Range("G145").Value = 0.051328403177
$ws.Range("E146").Value = 265728
$ws.Range("G146").Value = 0.052038161318
$ws.Range("E147").Value = 268750
$ws.Range("G147").Value = 0.050415917266
$ws.Range("E148").Value = 271804
$ws.Range("G148").Value = 0.048886382858
$ws.Range("E149").Value = 274811
$ws.Range("G149").Value = 0.047821406311
$ws.Range("E150").Value = 278128
$ws.Range("G150").Value = 0.049729972997
$ws.Range("E151").Value = 280950
$ws.Range("G151").Value = 0.047828069775
$ws.Range("E152").Value = 282401
$ws.Range("E153").Value = 283491
$ws.Range("E154").Value = 286679
$ws.Range("E155").Value = 290005
$ws.Range("E156").Value = 293036
$ws.Range("E157").Value = 296013
$ws.Range("E158").Value = 298882
$ws.Range("E159").Value = 300234
$ws.Range("E160").Value = 301335
$ws.Range("E161").Value = 302609
$ws.Range("E162").Value = 306951
$ws.Range("E163").Value = 310514
$ws.Range("E164").Value = 313836
$ws.Range("E165").Value = 316362
$ws.Range("E166").Value = 317785
$ws.Range("E167").Value = 318657
$ws.Range("E168").Value = 322097
$ws.Range("E169").Value = 326129
$ws.Range("E170").Value = 329492
$ws.Range("E171").Value = 333263
$ws.Range("E172").Value = 336569
$ws.Range("E173").Value = 338099
$ws.Range("F173").Value = 1530
$ws.Range("G173").Value = 0.031899182829
$ws.Range("E174").Value = 339264
$ws.Range("G174").Value = 0.031202989275
$ws.Range("E175").Value = 343671
$ws.Range("G175").Value = 0.03198294243
$ws.Range("E176").Value = 347489
$ws.Range("G176").Value = 0.030758426966
$ws.Range("E177").Value = 351041
$ws.Range("G177").Value = 0.030813494825
$ws.Range("E178").Value = 354808
$ws.Range("G178").Value = 0.030169412856
$ws.Range("E179").Value = 358152
$ws.Range("G179").Value = 0.029421303803
$ws.Range("E180").Value = 360045
$ws.Range("E181").Value = 361190
$ws.Range("E182").Value = 365320
$ws.Range("E183").Value = 369195
$ws.Range("E184").Value = 372798
$ws.Range("E185").Value = 376360
$ws.Range("F185").Value = 3562
$ws.Range("G185").Value = 0.022689309576
$ws.Range("E186").Value = 379938
$ws.Range("G186").Value = 0.021802992747
$ws.Range("E187").Value = 381945
$ws.Range("G187").Value = 0.020867579908
$ws.Range("E188").Value = 383346
$ws.Range("G188").Value = 0.021213215381
$ws.Range("E189").Value = 388190
$ws.Range("G189").Value = 0.018976825535
$ws.Range("E190").Value = 392883
$ws.Range("G190").Value = 0.017857142857
$ws.Range("E191").Value = 397107
$ws.Range("G191").Value = 0.018059155045
$ws.Range("E192").Value = 401423
$ws.Range("E193").Value = 405241
$ws.Range("E194").Value = 407053
$ws.Range("E195").Value = 408397
$ws.Range("E196").Value = 412970
$ws.Range("E197").Value = 417533
$ws.Range("E198").Value = 421816
$ws.Range("E199").Value = 425884
$ws.Range("E200").Value = 429931
$ws.Range("E201").Value = 431852
$ws.Range("E202").Value = 433249
$ws.Range("E203").Value = 438131
$ws.Range("E204").Value = 442976
$ws.Range("E205").Value = 447353
$ws.Range("E206").Value = 451330
$ws.Range("E207").Value = 455180
$ws.Range("E208").Value = 457150
$ws.Range("E209").Value = 458763
$ws.Range("E210").Value = 463723
$ws.Range("E211").Value = 468969
$ws.Range("F212").Value = 4593
$ws.Range("G212").Value = 0.021595635087
$ws.Range("G213").Value = 0.022509721806
$ws.Range("G214").Value = 0.021724938343
$ws.Range("G215").Value = 0.021599681078
$ws.Range("G216").Value = 0.022153380641
$ws.Range("G217").Value = 0.023122807017
$ws.Range("G218").Value = 0.024577194894
$ws.Range("E246").Value = 643316
$ws.Range("F246").Value = 8529
$ws.Range("G246").Value = 0.063275092322
$ws.Range("E247").Value = 651289
$ws.Range("G247").Value = 0.067231169109
$ws.Range("B248").Value = 27452
$ws.Range("C248").Value = 538
$ws.Range("D248").Value = 424.142857142857
$ws.Range("E248").Value = 658666
$ws.Range("G248").Value = 0.068587137312
$ws.Range("B249").Value = 28073
$ws.Range("D249").Value = 462.428571428571
$ws.Range("E249").Value = 666043
$ws.Range("G249").Value = 0.071211721224
$ws.Range("B250").Value = 28406
$ws.Range("D250").Value = 475.857142857142
$ws.Range("E250").Value = 669912
$ws.Range("G250").Value = 0.072816701278
$ws.Range("B251").Value = 28706
$ws.Range("D251").Value = 482.285714285714
$ws.Range("E251").Value = 673040
$ws.Range("G251").Value = 0.073113156469
$ws.Range("B252").Value = 29335
$ws.Range("D252").Value = 505.571428571428
$ws.Range("E252").Value = 681422
$ws.Range("G252").Value = 0.075887209177
$ws.Range("B253").Value = 29969
$ws.Range("C253").Value = 634
$ws.Range("D253").Value = 518.428571428571
$ws.Range("E253").Value = 690140
$ws.Range("G253").Value = 0.077502989919
$ws.Range("B254").Value = 30670
$ws.Range("D254").Value = 536.571428571428
$ws.Range("E254").Value = 698533
$ws.Range("G254").Value = 0.079502159004
$ws.Range("B255").Value = 31313
$ws.Range("C255").Value = 643
$ws.Range("D255").Value = 551.571428571428
$ws.Range("E255").Value = 706583
$ws.Range("F255").Value = 8050
$ws.Range("G255").Value = 0.080576830769
$ws.Range("B256").Value = 31997
$ws.Range("C256").Value = 684
$ws.Range("D256").Value = 560.571428571428
$ws.Range("E256").Value = 714254
$ws.Range("G256").Value = 0.081392213395
$ws.Range("B257").Value = 32394
$ws.Range("D257").Value = 569.714285714285
$ws.Range("E257").Value = 718083
$ws.Range("G257").Value = 0.08278839965900001
$ws.Range("B258").Value = 32670
$ws.Range("D258").Value = 566.285714285714
$ws.Range("E258").Value = 721215
$ws.Range("G258").Value = 0.08228334198200001
$ws.Range("B259").Value = 33415
$ws.Range("C259").Value = 745
$ws.Range("D259").Value = 582.857142857142
$ws.Range("E259").Value = 730177
$ws.Range("G259").Value = 0.083683724746
$ws.Range("B260").Value = 34139
$ws.Range("C260").Value = 724
$ws.Range("D260").Value = 595.714285714285
$ws.Range("E260").Value = 739509
$ws.Range("G260").Value = 0.08446596042
$ws.Range("B261").Value = 34807
$ws.Range("C261").Value = 668
$ws.Range("D261").Value = 591
$ws.Range("E261").Value = 748028
$ws.Range("G261").Value = 0.083584200424
$ws.Range("B262").Value = 35415
$ws.Range("D262").Value = 586
$ws.Range("E262").Value = 756742
$ws.Range("F262").Value = 8714
$ws.Range("G262").Value = 0.081779939791
$ws.Range("B263").Value = 35969
$ws.Range("C263").Value = 554
$ws.Range("D263").Value = 567.428571428571
$ws.Range("E263").Value = 764499
$ws.Range("G263").Value = 0.079052642053
$ws.Range("B264").Value = 36243
$ws.Range("D264").Value = 549.857142857142
$ws.Range("E264").Value = 768142
$ws.Range("G264").Value = 0.07688927065999999
$ws.Range("B265").Value = 36527
$ws.Range("D265").Value = 551
$ws.Range("E265").Value = 771791
$ws.Range("G265").Value = 0.076261467889
$ws.Range("B266").Value = 37145
$ws.Range("C266").Value = 618
$ws.Range("D266").Value = 532.857142857142
$ws.Range("E266").Value = 781581
$ws.Range("G266").Value = 0.072562446502
$ws.Range("B267").Value = 37899
$ws.Range("D267").Value = 537.142857142857
$ws.Range("E267").Value = 791078
$ws.Range("G267").Value = 0.072912020787
$ws.Range("B268").Value = 38490
$ws.Range("D268").Value = 526.142857142857
$ws.Range("E268").Value = 799371
$ws.Range("G268").Value = 0.07173324503799999
$ws.Range("B269").Value = 38893
$ws.Range("D269").Value = 496.857142857142
$ws.Range("E269").Value = 803962
$ws.Range("G269").Value = 0.073655230834
$ws.Range("B270").Value = 39047
$ws.Range("D270").Value = 439.714285714285
$ws.Range("E270").Value = 805013
$ws.Range("G270").Value = 0.075973737473
$ws.Range("B271").Value = 39467
$ws.Range("C271").Value = 420
$ws.Range("D271").Value = 460.571428571428
$ws.Range("E271").Value = 808536
$ws.Range("F271").Value = 3523
$ws.Range("G271").Value = 0.079813833737

# --- Append new row 272 (dimension will extend to A1:G272 automatically) ---
$ws.Range("A272").Value = 44192
$ws.Range("B272").Value = 39802
$ws.Range("C272").Value = 335
$ws.Range("D272").Value = 467.857142857142
$ws.Range("E272").Value = 811832
$ws.Range("F272").Value = 3296
$ws.Range("G272").Value = 0.081791164056
